$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.422.40'
$ws.Range("E2").Value = '  -4.16%  '
$ws.Range("D3").Value = '3.411.35'
$ws.Range("E3").Value = '  -3.79%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = '563.36'
$ws.Range("E5").Value = '  +0.73%  '
$ws.Range("D6").Value = '173.18'
$ws.Range("E6").Value = '  -8.95%  '
$ws.Range("D7").Value = '0.619'
$ws.Range("E7").Value = '  +0.52%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '0.621'
$ws.Range("E9").Value = '  -1.69%  '
$ws.Range("D10").Value = '0.155'
$ws.Range("E10").Value = '  +1.27%  '
$ws.Range("D11").Value = '57.07'
$ws.Range("E11").Value = '  +3.83%  '
$ws.Range("D12").Value = '''0.0000270'
$ws.Range("E12").Value = '  -0.88%  '
$ws.Range("D13").Value = '9.04'
$ws.Range("E13").Value = '  -3.66%  '
$ws.Range("D14").Value = '3.963.25'
$ws.Range("E14").Value = '  -3.38%  '
$ws.Range("E15").Value = '  -1.18%  '
$ws.Range("D16").Value = '3.416.73'
$ws.Range("E16").Value = '  -3.47%  '
$ws.Range("D17").Value = '17.99'
$ws.Range("E17").Value = '  -1.46%  '
$ws.Range("D18").Value = '11.81'
$ws.Range("E18").Value = '  -2.20%  '
$ws.Range("D19").Value = '64.448.87'
$ws.Range("E19").Value = '  -4.12%  '
$ws.Range("D20").Value = '0.987'
$ws.Range("E20").Value = '  -1.18%  '
$ws.Range("D21").Value = '408.74'
$ws.Range("E21").Value = '  -4.87%  '
$ws.Range("D22").Value = '4.14'
$ws.Range("E22").Value = '  +0.62%  '
$ws.Range("D23").Value = '4.39'
$ws.Range("E23").Value = '  +5.94%  '
$ws.Range("D24").Value = '13.37'
$ws.Range("E24").Value = '  +8.00%  '
$ws.Range("D25").Value = '83.07'
$ws.Range("E25").Value = '  -2.47%  '
$ws.Range("D26").Value = '10.73'
$ws.Range("E26").Value = '  -3.07%  '
$ws.Range("D27").Value = '2.76'
$ws.Range("E27").Value = '  -4.71%  '
$ws.Range("D28").Value = '8.85'
$ws.Range("E28").Value = '  -2.01%  '
$ws.Range("D29").Value = '29.61'
$ws.Range("E29").Value = '  -2.91%  '
$ws.Range("D30").Value = '6.68'
$ws.Range("E30").Value = '  +1.35%  '
$ws.Range("D31").Value = '''590.60'
$ws.Range("E31").Value = '  -8.25%  '
$ws.Range("D32").Value = '11.47'
$ws.Range("E32").Value = '  -2.20%  '
$ws.Range("D33").Value = '0.107'
$ws.Range("E33").Value = '  -3.60%  '
$ws.Range("E34").Value = '  +5.20%  '
$ws.Range("D35").Value = '58.95'
$ws.Range("E35").Value = '  -1.95%  '
$ws.Range("D36").Value = '''1.00'
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("D37").Value = '35.75'
$ws.Range("E37").Value = '  -6.74%  '
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").Value = '3.42'
$ws.Range("E38").Value = '  +1.28%  '
$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").Value = '0.372'
$ws.Range("E39").Value = '  -4.39%  '
$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0745'
$ws.Range("E40").Value = '  -8.44%  '
$ws.Range("D41").Value = '3.165.00'
$ws.Range("E41").Value = '  +1.49%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("D43").Value = '''2.90'
$ws.Range("E43").Value = '  +0.97%  '
$ws.Range("D44").Value = '2.53'
$ws.Range("E44").Value = '  -4.74%  '
$ws.Range("D45").Value = '3.24'
$ws.Range("E45").Value = '  -4.16%  '
$ws.Range("D46").Value = '0.0407'
$ws.Range("E46").Value = '  -2.69%  '
$ws.Range("E47").Value = '  -4.56%  '
$ws.Range("D48").Value = '0.129'
$ws.Range("E48").Value = '  -1.53%  '
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").Value = '8.31'
$ws.Range("E49").Value = '  -3.59%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '135.25'
$ws.Range("E50").Value = '  -4.28%  '
$ws.Range("D51").Value = '2.76'
$ws.Range("E51").Value = '  +1.52%  '
